$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, shifting existing rows 141-143 down to 142-144,
# and carrying the formatting (date style on column D) down with them.
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with its data. Columns A,B,C,E,F,G,H,I,
# N,Q,R are identical to the surrounding rows for this series.
$ws.Range("A141").Value = 5
$ws.Range("B141").Value = "Macroferia Regional de Talca"
$ws.Range("C141").Value = "Maule"
$ws.Range("D141").Value = 44595
$ws.Range("E141").Value = 7
$ws.Range("F141").Value = 100112031
$ws.Range("G141").Value = "Poroto verde"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 100
$ws.Range("K141").Value = 33000
$ws.Range("L141").Value = 33000
$ws.Range("M141").Value = 33000
$ws.Range("N141").Value = "$/saco 25 kilos"
$ws.Range("O141").Value = "Región del Maule"
$ws.Range("P141").Value = 1320
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"
